$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header from "percent_recovery" to "enrichment"
$ws.Range("C1").Value = "enrichment"

# Divide each numeric value in C2:C11 by 100 (percent_recovery -> enrichment fraction)
for ($r = 2; $r -le 11; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $cell.Value2 = $cell.Value2 / 100
}
